$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns (row 1) to the new machine-friendly names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the Spanish connector words ("de", "del", "la", "las", "los",
#    "el", "y") inside the state/municipality text columns (A and B) for the
#    data rows. Numeric columns (C, D) are left untouched.
for ($r = 2; $r -le 684; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $orig = $cell.Text
        if ($orig -ne "") {
            $new = $orig -replace '\bde\b','De'
            $new = $new -replace '\bdel\b','Del'
            $new = $new -replace '\blas\b','Las'
            $new = $new -replace '\bla\b','La'
            $new = $new -replace '\blos\b','Los'
            $new = $new -replace '\bel\b','El'
            $new = $new -replace '\by\b','Y'
            $unchanged = $orig.Equals($new)
            if ($unchanged -eq $false) {
                $cell.Value = $new
            }
        }
    }
}

# 3) Drop the trailing footnote rows (sample size / source / author / date)
#    that used to sit below the data table, shrinking the sheet back down to
#    the real data range.
$ws.Range("A685:A690").EntireRow.Delete()
